# Configurate istanze e librerire ISPRO
#
# 1. Add a new "ISPRO" row to the "r CustomerUnit_AnalysisUnit" sheet,
#    mirroring the existing BE-COUNTERPARTY_* rows.
# 2. Switch the active / selected worksheet from
#    "r CustomerUnit_AnalysisUnit" back to "Customer_Unit".
# 3. Update the remembered cell selection on both worksheets.

$wb = $excel.ActiveWorkbook

$wsCustomerUnit = $wb.Worksheets.Item("Customer_Unit")
$wsAnalysisUnit = $wb.Worksheets.Item("r CustomerUnit_AnalysisUnit")

# --- 1. New data row (row 9) on the relation sheet -----------------------
$wsAnalysisUnit.Range("A9").Value = "CREATE/MODIFY"
$wsAnalysisUnit.Range("B9").Value = "BE-COUNTERPARTY_ISPRO"
$wsAnalysisUnit.Range("C9").Value = "BE-COUNTERPARTY_ISPRO"
$wsAnalysisUnit.Range("D9").Value = "BE-COUNTERPARTY_ISPRO"
$wsAnalysisUnit.Range("E9").Value = "BE"
$wsAnalysisUnit.Range("F9").Value = "COUNTERPARTY_ISPRO"

# --- 2 & 3. Activate sheets / restore per-sheet selection -----------------
# Touch the analysis-unit sheet's selection first (it currently holds the
# active tab), then finish on Customer_Unit so it becomes the active tab,
# matching the saved view state in the workbook.
[void]$wsAnalysisUnit.Activate()
[void]$wsAnalysisUnit.Range("E13").Select()

[void]$wsCustomerUnit.Activate()
[void]$wsCustomerUnit.Range("D15").Select()
